$d = $word.ActiveDocument
$nbsp = [char]160

# --- Article 9 title: "Deux cours collectifs" -> "Cours collectifs" ---
$oldTitle = "Article 9" + $nbsp + ": Deux cours collectifs offerts"
$newTitle = "Article 9" + $nbsp + ": Cours collectifs offerts"
$d.Content.Find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)

# --- Body paragraph describing the offered group lessons ---
$oldBody = "Une fois le programme terminé vous disposez de 30 jours pour effectuer 2 cours collectifs qui vous sont offerts. Ces deux cours doivent impérativement être consécutifs au programme (de date à date). Ils ne peuvent en aucun cas être repoussés à une date ultérieure.  "
$newBody = "Une fois le programme terminé, 4 cours collectifs qui vous sont offerts. Ces cours doivent impérativement être consécutifs au programme. Vous disposez de deux mois pour les suivre. Ils ne peuvent en aucun cas être repoussés à une date ultérieure."
$d.Content.Find.Execute($oldBody, $true, $false, $false, $false, $false, $true, 1, $false, $newBody, 2)
